# MHD2-259: Report template and related changes for reporting on 136 genes
#
# The clinical-context header table used an amber/theme-tinted shading
# (theme accent4 @ 33% tint => FFF2CC) for the table's default shading
# and a slightly different grey (E8E7EC) for its title cell. Both are
# being replaced with the new template's flat lavender-grey fill
# (ECEAF2) so the table background and the title cell match.

$d = $word.ActiveDocument

# New fill colour: ECEAF2, packed as a Word/VBA OLE colour (0x00BBGGRR).
$newFill = 15919852   # RGB(0xEC, 0xEA, 0xF2)

$tbl = $d.Tables.Item(1)

# Table-wide default shading (writes <w:tblPr><w:shd .../>).
$tbl.Rows.Shading.BackgroundPatternColor = $newFill

# Shading of the (single) header cell, which carries its own explicit
# override (writes <w:tcPr><w:shd .../> on the first cell).
$cell = $tbl.Cell(1, 1)
$cell.Shading.BackgroundPatternColor = $newFill
